$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price cells we will rewrite, so that
# values like "65.50" or "12.17" stay exact text instead of being
# auto-converted to numbers by Excel.
$priceCells = @("D2", "D3", "D5", "D6", "D10", "D11", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D39", "D41", "D44", "D45", "D47", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row, in the same order as the source diff.
$ws.Range("D2").Value = '44.283.09'
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").Value = '2.241.36'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = '306.93'
$ws.Range("E5").Value = '  -2.98%  '

$ws.Range("D6").Value = '94.91'
$ws.Range("E6").Value = '  -4.62%  '

$ws.Range("E7").Value = '  -0.78%  '

$ws.Range("E8").Value = '  +0.26%  '

$ws.Range("E9").Value = '  -1.83%  '

$ws.Range("D10").Value = '34.63'
$ws.Range("E10").Value = '  -4.74%  '

$ws.Range("D11").Value = '0.0811'
$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("E12").Value = '  -2.85%  '

$ws.Range("E13").Value = '  -0.11%  '

$ws.Range("D14").Value = '2.337.11'
$ws.Range("E14").Value = '  +3.71%  '

$ws.Range("D15").Value = '2.583.33'
$ws.Range("E15").Value = '  -0.15%  '

$ws.Range("D16").Value = '0.829'
$ws.Range("E16").Value = '  -2.10%  '

$ws.Range("E17").Value = '  -3.77%  '

$ws.Range("D18").Value = '44.023.66'
$ws.Range("E18").Value = '  +0.12%  '

$ws.Range("E19").Value = '  -1.55%  '

$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").Value = '12.17'
$ws.Range("E20").Value = '  -7.64%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.38'
$ws.Range("E21").Value = '  +0.64%  '

$ws.Range("D22").Value = '65.50'

$ws.Range("D23").Value = '237.65'
$ws.Range("E23").Value = '  -0.28%  '

$ws.Range("D24").Value = '2.94'
$ws.Range("E24").Value = '  -1.68%  '

$ws.Range("D25").Value = '2.01'
$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").Value = '39.57'
$ws.Range("E27").Value = '  +8.31%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  +3.89%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '9.91'
$ws.Range("E29").Value = '  -2.57%  '

$ws.Range("D30").Value = '20.02'
$ws.Range("E30").Value = '  -0.40%  '

$ws.Range("D31").Value = '5.85'
$ws.Range("E31").Value = '  -2.69%  '

$ws.Range("D32").Value = '153.08'
$ws.Range("E32").Value = '  -2.15%  '

$ws.Range("D33").Value = '0.0796'
$ws.Range("E33").Value = '  -5.77%  '

$ws.Range("E34").Value = '  -1.88%  '

$ws.Range("D35").Value = '3.13'
$ws.Range("E35").Value = '  -5.09%  '

$ws.Range("E36").Value = '  +1.63%  '

$ws.Range("E37").Value = '  -0.32%  '

$ws.Range("E38").Value = '  -7.58%  '

$ws.Range("D39").Value = '3.51'
$ws.Range("E39").Value = '  -1.15%  '

$ws.Range("E40").Value = '  -5.10%  '

$ws.Range("D41").Value = '14.32'
$ws.Range("E41").Value = '  -7.35%  '

$ws.Range("E42").Value = '  -3.44%  '

$ws.Range("E43").Value = '  +0.31%  '

$ws.Range("D44").Value = '1.744.30'
$ws.Range("E44").Value = '  +2.14%  '

$ws.Range("D45").Value = '83.09'
$ws.Range("E45").Value = '  +0.19%  '

$ws.Range("E46").Value = '  -1.86%  '

$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").Value = '4.93'
$ws.Range("E47").Value = '  -5.08%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '99.51'
$ws.Range("E48").Value = '  -2.57%  '

$ws.Range("E49").Value = '  -0.66%  '

$ws.Range("D50").Value = '54.84'
$ws.Range("E50").Value = '  -3.46%  '

$ws.Range("D51").Value = '8.07'
$ws.Range("E51").Value = '  -0.78%  '

